# Updated cryptos list on Thu Jun 20 11:45:14 UTC 2024 with GitHub Actions
#
# Refreshes the live price / 1h-volume columns (D, E) for every coin row on
# Sheet1, and fixes the rank-43/44 Coin+Link+Price+Volume rows which had
# swapped (EnergySwap <-> ONDO), plus replaces row 51 (InjectiveProtocol)
# with VeChain's current data.
#
# Cells are plain text in the source workbook (t="inlineStr"), e.g. Price
# "3.608.00" uses dots as thousands separators, so it never parses as a
# genuine number - but some refreshed prices (e.g. "605.24", "0.0274") DO
# look like valid numbers. Excel's Range.Value setter auto-converts such
# strings to the Number type, which would silently change the stored cell
# type/format. To keep those cells as text (matching the original file),
# values that parse as numbers are entered with a leading apostrophe (the
# standard Excel "treat as text" quote-prefix), then ClearFormats() removes
# the transient quote-prefix cell style so formatting stays identical to
# the untouched cells around it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.243.55'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").Value = '3.608.00'
$ws.Range("E3").Value = '  +2.41%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''605.24'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("D6").Value = '''139.72'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.75%  '
$ws.Range("D7").Value = '3.612.32'
$ws.Range("E7").Value = '  +2.53%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '''0.501'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("D10").Value = '''0.127'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.04%  '
$ws.Range("D11").Value = '''7.23'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.19%  '
$ws.Range("E12").Value = '  +2.64%  '
$ws.Range("D13").Value = '4.219.58'
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("D14").Value = '''28.52'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.41%  '
$ws.Range("E15").Value = '  +3.37%  '
$ws.Range("D16").Value = '3.606.44'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '66.354.97'
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("D19").Value = '''10.18'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("D20").Value = '''14.68'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.42%  '
$ws.Range("D21").Value = '''5.93'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("D22").Value = '''398.55'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.02%  '
$ws.Range("D23").Value = '''0.592'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.62%  '
$ws.Range("D24").Value = '3.754.40'
$ws.Range("E24").Value = '  +2.34%  '
$ws.Range("D25").Value = '''75.16'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +6.43%  '
$ws.Range("D28").Value = '''8.19'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.38%  '
$ws.Range("E29").Value = '  +28.30%  '
$ws.Range("D30").Value = '''8.69'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +7.11%  '
$ws.Range("D31").Value = '''2.35'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.58%  '
$ws.Range("D32").Value = '''0.999'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("D33").Value = '3.613.48'
$ws.Range("E33").Value = '  +1.91%  '
$ws.Range("D34").Value = '''24.69'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.81%  '
$ws.Range("E35").Value = '  +4.87%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").Value = '''5.42'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +9.12%  '
$ws.Range("D38").Value = '''1.64'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +5.52%  '
$ws.Range("D39").Value = '''7.08'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("D40").Value = '''168.77'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("E41").Value = '  +5.83%  '
$ws.Range("D42").Value = '''0.846'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.82%  '
$ws.Range("B43").Value = 'ONDO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D43").Value = '''1.28'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +7.61%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''26.21'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").Value = '''43.24'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("D46").Value = '''4.58'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.84%  '
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("E48").Value = '  +3.97%  '
$ws.Range("D49").Value = '''7.05'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.04%  '
$ws.Range("D50").Value = '2.465.57'
$ws.Range("E50").Value = '  +3.32%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '''0.0274'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.38%  '
